$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2023-11-04 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-05 Sunday", 2) | Out-Null

# Update the 100 math-problem cells in the table, in document order
$t = $d.Tables.Item(1)
$mismatches = 0

$cell = $t.Cell(1, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "67-27=40") { $mismatches++; Write-Output "MISMATCH at (1,1): expected '67-27=40' got '$cellText'" }
$cell.Range.Text = "89-76=13"

$cell = $t.Cell(1, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "0+22=22") { $mismatches++; Write-Output "MISMATCH at (1,2): expected '0+22=22' got '$cellText'" }
$cell.Range.Text = "88-4=84"

$cell = $t.Cell(1, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "30+2=32") { $mismatches++; Write-Output "MISMATCH at (1,3): expected '30+2=32' got '$cellText'" }
$cell.Range.Text = "34+23=57"

$cell = $t.Cell(1, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "11+48=59") { $mismatches++; Write-Output "MISMATCH at (1,4): expected '11+48=59' got '$cellText'" }
$cell.Range.Text = "28+65=93"

$cell = $t.Cell(1, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "6+76=82") { $mismatches++; Write-Output "MISMATCH at (1,5): expected '6+76=82' got '$cellText'" }
$cell.Range.Text = "28+54=82"

$cell = $t.Cell(2, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "65+11=76") { $mismatches++; Write-Output "MISMATCH at (2,1): expected '65+11=76' got '$cellText'" }
$cell.Range.Text = "61+18=79"

$cell = $t.Cell(2, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "11+31=42") { $mismatches++; Write-Output "MISMATCH at (2,2): expected '11+31=42' got '$cellText'" }
$cell.Range.Text = "48-7=41"

$cell = $t.Cell(2, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "54-36=18") { $mismatches++; Write-Output "MISMATCH at (2,3): expected '54-36=18' got '$cellText'" }
$cell.Range.Text = "79+13=92"

$cell = $t.Cell(2, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "95-47=48") { $mismatches++; Write-Output "MISMATCH at (2,4): expected '95-47=48' got '$cellText'" }
$cell.Range.Text = "70-39=31"

$cell = $t.Cell(2, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "84-52=32") { $mismatches++; Write-Output "MISMATCH at (2,5): expected '84-52=32' got '$cellText'" }
$cell.Range.Text = "16+3=19"

$cell = $t.Cell(3, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "11+13=24") { $mismatches++; Write-Output "MISMATCH at (3,1): expected '11+13=24' got '$cellText'" }
$cell.Range.Text = "73-34=39"

$cell = $t.Cell(3, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "22-3=19") { $mismatches++; Write-Output "MISMATCH at (3,2): expected '22-3=19' got '$cellText'" }
$cell.Range.Text = "34+14=48"

$cell = $t.Cell(3, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "94-47=47") { $mismatches++; Write-Output "MISMATCH at (3,3): expected '94-47=47' got '$cellText'" }
$cell.Range.Text = "22+6=28"

$cell = $t.Cell(3, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "72-29=43") { $mismatches++; Write-Output "MISMATCH at (3,4): expected '72-29=43' got '$cellText'" }
$cell.Range.Text = "1+19=20"

$cell = $t.Cell(3, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "64+27=91") { $mismatches++; Write-Output "MISMATCH at (3,5): expected '64+27=91' got '$cellText'" }
$cell.Range.Text = "66-20=46"

$cell = $t.Cell(4, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "95-78=17") { $mismatches++; Write-Output "MISMATCH at (4,1): expected '95-78=17' got '$cellText'" }
$cell.Range.Text = "40+31=71"

$cell = $t.Cell(4, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "42-29=13") { $mismatches++; Write-Output "MISMATCH at (4,2): expected '42-29=13' got '$cellText'" }
$cell.Range.Text = "22-15=7"

$cell = $t.Cell(4, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "25+68=93") { $mismatches++; Write-Output "MISMATCH at (4,3): expected '25+68=93' got '$cellText'" }
$cell.Range.Text = "49-34=15"

$cell = $t.Cell(4, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "43+56=99") { $mismatches++; Write-Output "MISMATCH at (4,4): expected '43+56=99' got '$cellText'" }
$cell.Range.Text = "47+29=76"

$cell = $t.Cell(4, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "23+47=70") { $mismatches++; Write-Output "MISMATCH at (4,5): expected '23+47=70' got '$cellText'" }
$cell.Range.Text = "95-63=32"

$cell = $t.Cell(5, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "94-26=68") { $mismatches++; Write-Output "MISMATCH at (5,1): expected '94-26=68' got '$cellText'" }
$cell.Range.Text = "20+65=85"

$cell = $t.Cell(5, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "56-18=38") { $mismatches++; Write-Output "MISMATCH at (5,2): expected '56-18=38' got '$cellText'" }
$cell.Range.Text = "67-24=43"

$cell = $t.Cell(5, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "15+76=91") { $mismatches++; Write-Output "MISMATCH at (5,3): expected '15+76=91' got '$cellText'" }
$cell.Range.Text = "71-61=10"

$cell = $t.Cell(5, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "88-58=30") { $mismatches++; Write-Output "MISMATCH at (5,4): expected '88-58=30' got '$cellText'" }
$cell.Range.Text = "75-46=29"

$cell = $t.Cell(5, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "48-14=34") { $mismatches++; Write-Output "MISMATCH at (5,5): expected '48-14=34' got '$cellText'" }
$cell.Range.Text = "23+12=35"

$cell = $t.Cell(6, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "79+17=96") { $mismatches++; Write-Output "MISMATCH at (6,1): expected '79+17=96' got '$cellText'" }
$cell.Range.Text = "96-64=32"

$cell = $t.Cell(6, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "88-34=54") { $mismatches++; Write-Output "MISMATCH at (6,2): expected '88-34=54' got '$cellText'" }
$cell.Range.Text = "75-33=42"

$cell = $t.Cell(6, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "64-63=1") { $mismatches++; Write-Output "MISMATCH at (6,3): expected '64-63=1' got '$cellText'" }
$cell.Range.Text = "17+3=20"

$cell = $t.Cell(6, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "67+10=77") { $mismatches++; Write-Output "MISMATCH at (6,4): expected '67+10=77' got '$cellText'" }
$cell.Range.Text = "20-5=15"

$cell = $t.Cell(6, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "37-5=32") { $mismatches++; Write-Output "MISMATCH at (6,5): expected '37-5=32' got '$cellText'" }
$cell.Range.Text = "76+15=91"

$cell = $t.Cell(7, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "6+17=23") { $mismatches++; Write-Output "MISMATCH at (7,1): expected '6+17=23' got '$cellText'" }
$cell.Range.Text = "74+17=91"

$cell = $t.Cell(7, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "3+12=15") { $mismatches++; Write-Output "MISMATCH at (7,2): expected '3+12=15' got '$cellText'" }
$cell.Range.Text = "58-54=4"

$cell = $t.Cell(7, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "60-11=49") { $mismatches++; Write-Output "MISMATCH at (7,3): expected '60-11=49' got '$cellText'" }
$cell.Range.Text = "49+38=87"

$cell = $t.Cell(7, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "21+38=59") { $mismatches++; Write-Output "MISMATCH at (7,4): expected '21+38=59' got '$cellText'" }
$cell.Range.Text = "97-74=23"

$cell = $t.Cell(7, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "37+41=78") { $mismatches++; Write-Output "MISMATCH at (7,5): expected '37+41=78' got '$cellText'" }
$cell.Range.Text = "33-20=13"

$cell = $t.Cell(8, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "79-11=68") { $mismatches++; Write-Output "MISMATCH at (8,1): expected '79-11=68' got '$cellText'" }
$cell.Range.Text = "80-14=66"

$cell = $t.Cell(8, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "27-0=27") { $mismatches++; Write-Output "MISMATCH at (8,2): expected '27-0=27' got '$cellText'" }
$cell.Range.Text = "63-39=24"

$cell = $t.Cell(8, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "40-27=13") { $mismatches++; Write-Output "MISMATCH at (8,3): expected '40-27=13' got '$cellText'" }
$cell.Range.Text = "96-60=36"

$cell = $t.Cell(8, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "48-42=6") { $mismatches++; Write-Output "MISMATCH at (8,4): expected '48-42=6' got '$cellText'" }
$cell.Range.Text = "70-2=68"

$cell = $t.Cell(8, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "28+62=90") { $mismatches++; Write-Output "MISMATCH at (8,5): expected '28+62=90' got '$cellText'" }
$cell.Range.Text = "39+27=66"

$cell = $t.Cell(9, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "25-7=18") { $mismatches++; Write-Output "MISMATCH at (9,1): expected '25-7=18' got '$cellText'" }
$cell.Range.Text = "28-12=16"

$cell = $t.Cell(9, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "60-43=17") { $mismatches++; Write-Output "MISMATCH at (9,2): expected '60-43=17' got '$cellText'" }
$cell.Range.Text = "58+17=75"

$cell = $t.Cell(9, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "18+70=88") { $mismatches++; Write-Output "MISMATCH at (9,3): expected '18+70=88' got '$cellText'" }
$cell.Range.Text = "32+23=55"

$cell = $t.Cell(9, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "83-1=82") { $mismatches++; Write-Output "MISMATCH at (9,4): expected '83-1=82' got '$cellText'" }
$cell.Range.Text = "68-50=18"

$cell = $t.Cell(9, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "40+7=47") { $mismatches++; Write-Output "MISMATCH at (9,5): expected '40+7=47' got '$cellText'" }
$cell.Range.Text = "92-46=46"

$cell = $t.Cell(10, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "71-7=64") { $mismatches++; Write-Output "MISMATCH at (10,1): expected '71-7=64' got '$cellText'" }
$cell.Range.Text = "55+30=85"

$cell = $t.Cell(10, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "64+20=84") { $mismatches++; Write-Output "MISMATCH at (10,2): expected '64+20=84' got '$cellText'" }
$cell.Range.Text = "42+31=73"

$cell = $t.Cell(10, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "57+18=75") { $mismatches++; Write-Output "MISMATCH at (10,3): expected '57+18=75' got '$cellText'" }
$cell.Range.Text = "42+55=97"

$cell = $t.Cell(10, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "32+48=80") { $mismatches++; Write-Output "MISMATCH at (10,4): expected '32+48=80' got '$cellText'" }
$cell.Range.Text = "80+0=80"

$cell = $t.Cell(10, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "83+2=85") { $mismatches++; Write-Output "MISMATCH at (10,5): expected '83+2=85' got '$cellText'" }
$cell.Range.Text = "38-4=34"

$cell = $t.Cell(11, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "98-6=92") { $mismatches++; Write-Output "MISMATCH at (11,1): expected '98-6=92' got '$cellText'" }
$cell.Range.Text = "43+25=68"

$cell = $t.Cell(11, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "80+7=87") { $mismatches++; Write-Output "MISMATCH at (11,2): expected '80+7=87' got '$cellText'" }
$cell.Range.Text = "2+49=51"

$cell = $t.Cell(11, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "69-69=0") { $mismatches++; Write-Output "MISMATCH at (11,3): expected '69-69=0' got '$cellText'" }
$cell.Range.Text = "33-24=9"

$cell = $t.Cell(11, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "54-23=31") { $mismatches++; Write-Output "MISMATCH at (11,4): expected '54-23=31' got '$cellText'" }
$cell.Range.Text = "68-33=35"

$cell = $t.Cell(11, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "31+21=52") { $mismatches++; Write-Output "MISMATCH at (11,5): expected '31+21=52' got '$cellText'" }
$cell.Range.Text = "38+39=77"

$cell = $t.Cell(12, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "74-20=54") { $mismatches++; Write-Output "MISMATCH at (12,1): expected '74-20=54' got '$cellText'" }
$cell.Range.Text = "96-18=78"

$cell = $t.Cell(12, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "36+57=93") { $mismatches++; Write-Output "MISMATCH at (12,2): expected '36+57=93' got '$cellText'" }
$cell.Range.Text = "80-79=1"

$cell = $t.Cell(12, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "76-53=23") { $mismatches++; Write-Output "MISMATCH at (12,3): expected '76-53=23' got '$cellText'" }
$cell.Range.Text = "60-40=20"

$cell = $t.Cell(12, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "53+5=58") { $mismatches++; Write-Output "MISMATCH at (12,4): expected '53+5=58' got '$cellText'" }
$cell.Range.Text = "6+81=87"

$cell = $t.Cell(12, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "82-64=18") { $mismatches++; Write-Output "MISMATCH at (12,5): expected '82-64=18' got '$cellText'" }
$cell.Range.Text = "33+42=75"

$cell = $t.Cell(13, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "32+36=68") { $mismatches++; Write-Output "MISMATCH at (13,1): expected '32+36=68' got '$cellText'" }
$cell.Range.Text = "38+3=41"

$cell = $t.Cell(13, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "83-76=7") { $mismatches++; Write-Output "MISMATCH at (13,2): expected '83-76=7' got '$cellText'" }
$cell.Range.Text = "52-46=6"

$cell = $t.Cell(13, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "92-8=84") { $mismatches++; Write-Output "MISMATCH at (13,3): expected '92-8=84' got '$cellText'" }
$cell.Range.Text = "63-39=24"

$cell = $t.Cell(13, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "13+24=37") { $mismatches++; Write-Output "MISMATCH at (13,4): expected '13+24=37' got '$cellText'" }
$cell.Range.Text = "59+29=88"

$cell = $t.Cell(13, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "71-35=36") { $mismatches++; Write-Output "MISMATCH at (13,5): expected '71-35=36' got '$cellText'" }
$cell.Range.Text = "70+6=76"

$cell = $t.Cell(14, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "66-38=28") { $mismatches++; Write-Output "MISMATCH at (14,1): expected '66-38=28' got '$cellText'" }
$cell.Range.Text = "87+5=92"

$cell = $t.Cell(14, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "46-34=12") { $mismatches++; Write-Output "MISMATCH at (14,2): expected '46-34=12' got '$cellText'" }
$cell.Range.Text = "39-24=15"

$cell = $t.Cell(14, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "57-46=11") { $mismatches++; Write-Output "MISMATCH at (14,3): expected '57-46=11' got '$cellText'" }
$cell.Range.Text = "36+28=64"

$cell = $t.Cell(14, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "57+5=62") { $mismatches++; Write-Output "MISMATCH at (14,4): expected '57+5=62' got '$cellText'" }
$cell.Range.Text = "13+15=28"

$cell = $t.Cell(14, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "64-54=10") { $mismatches++; Write-Output "MISMATCH at (14,5): expected '64-54=10' got '$cellText'" }
$cell.Range.Text = "22+47=69"

$cell = $t.Cell(15, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "3+5=8") { $mismatches++; Write-Output "MISMATCH at (15,1): expected '3+5=8' got '$cellText'" }
$cell.Range.Text = "66+32=98"

$cell = $t.Cell(15, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "45-25=20") { $mismatches++; Write-Output "MISMATCH at (15,2): expected '45-25=20' got '$cellText'" }
$cell.Range.Text = "93-77=16"

$cell = $t.Cell(15, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "86-16=70") { $mismatches++; Write-Output "MISMATCH at (15,3): expected '86-16=70' got '$cellText'" }
$cell.Range.Text = "4+87=91"

$cell = $t.Cell(15, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "74-9=65") { $mismatches++; Write-Output "MISMATCH at (15,4): expected '74-9=65' got '$cellText'" }
$cell.Range.Text = "11+55=66"

$cell = $t.Cell(15, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "28+56=84") { $mismatches++; Write-Output "MISMATCH at (15,5): expected '28+56=84' got '$cellText'" }
$cell.Range.Text = "73-32=41"

$cell = $t.Cell(16, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "10+66=76") { $mismatches++; Write-Output "MISMATCH at (16,1): expected '10+66=76' got '$cellText'" }
$cell.Range.Text = "29-10=19"

$cell = $t.Cell(16, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "48-19=29") { $mismatches++; Write-Output "MISMATCH at (16,2): expected '48-19=29' got '$cellText'" }
$cell.Range.Text = "99-67=32"

$cell = $t.Cell(16, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "11+52=63") { $mismatches++; Write-Output "MISMATCH at (16,3): expected '11+52=63' got '$cellText'" }
$cell.Range.Text = "39+46=85"

$cell = $t.Cell(16, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "20-16=4") { $mismatches++; Write-Output "MISMATCH at (16,4): expected '20-16=4' got '$cellText'" }
$cell.Range.Text = "14+31=45"

$cell = $t.Cell(16, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "23+72=95") { $mismatches++; Write-Output "MISMATCH at (16,5): expected '23+72=95' got '$cellText'" }
$cell.Range.Text = "81-36=45"

$cell = $t.Cell(17, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "36+6=42") { $mismatches++; Write-Output "MISMATCH at (17,1): expected '36+6=42' got '$cellText'" }
$cell.Range.Text = "36+48=84"

$cell = $t.Cell(17, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "18-2=16") { $mismatches++; Write-Output "MISMATCH at (17,2): expected '18-2=16' got '$cellText'" }
$cell.Range.Text = "77-11=66"

$cell = $t.Cell(17, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "65-1=64") { $mismatches++; Write-Output "MISMATCH at (17,3): expected '65-1=64' got '$cellText'" }
$cell.Range.Text = "29+11=40"

$cell = $t.Cell(17, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "60-50=10") { $mismatches++; Write-Output "MISMATCH at (17,4): expected '60-50=10' got '$cellText'" }
$cell.Range.Text = "92-41=51"

$cell = $t.Cell(17, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "62+12=74") { $mismatches++; Write-Output "MISMATCH at (17,5): expected '62+12=74' got '$cellText'" }
$cell.Range.Text = "91+5=96"

$cell = $t.Cell(18, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "71-39=32") { $mismatches++; Write-Output "MISMATCH at (18,1): expected '71-39=32' got '$cellText'" }
$cell.Range.Text = "95-84=11"

$cell = $t.Cell(18, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "82-33=49") { $mismatches++; Write-Output "MISMATCH at (18,2): expected '82-33=49' got '$cellText'" }
$cell.Range.Text = "71-42=29"

$cell = $t.Cell(18, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "70-55=15") { $mismatches++; Write-Output "MISMATCH at (18,3): expected '70-55=15' got '$cellText'" }
$cell.Range.Text = "19-0=19"

$cell = $t.Cell(18, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "77-13=64") { $mismatches++; Write-Output "MISMATCH at (18,4): expected '77-13=64' got '$cellText'" }
$cell.Range.Text = "37+62=99"

$cell = $t.Cell(18, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "10+36=46") { $mismatches++; Write-Output "MISMATCH at (18,5): expected '10+36=46' got '$cellText'" }
$cell.Range.Text = "11+70=81"

$cell = $t.Cell(19, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "31-11=20") { $mismatches++; Write-Output "MISMATCH at (19,1): expected '31-11=20' got '$cellText'" }
$cell.Range.Text = "36+0=36"

$cell = $t.Cell(19, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "39+26=65") { $mismatches++; Write-Output "MISMATCH at (19,2): expected '39+26=65' got '$cellText'" }
$cell.Range.Text = "12+21=33"

$cell = $t.Cell(19, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "94-56=38") { $mismatches++; Write-Output "MISMATCH at (19,3): expected '94-56=38' got '$cellText'" }
$cell.Range.Text = "72-31=41"

$cell = $t.Cell(19, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "25+73=98") { $mismatches++; Write-Output "MISMATCH at (19,4): expected '25+73=98' got '$cellText'" }
$cell.Range.Text = "67-63=4"

$cell = $t.Cell(19, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "60+37=97") { $mismatches++; Write-Output "MISMATCH at (19,5): expected '60+37=97' got '$cellText'" }
$cell.Range.Text = "17+1=18"

$cell = $t.Cell(20, 1)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "98-6=92") { $mismatches++; Write-Output "MISMATCH at (20,1): expected '98-6=92' got '$cellText'" }
$cell.Range.Text = "24+66=90"

$cell = $t.Cell(20, 2)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "32+43=75") { $mismatches++; Write-Output "MISMATCH at (20,2): expected '32+43=75' got '$cellText'" }
$cell.Range.Text = "96-48=48"

$cell = $t.Cell(20, 3)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "59-22=37") { $mismatches++; Write-Output "MISMATCH at (20,3): expected '59-22=37' got '$cellText'" }
$cell.Range.Text = "48-4=44"

$cell = $t.Cell(20, 4)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "75-68=7") { $mismatches++; Write-Output "MISMATCH at (20,4): expected '75-68=7' got '$cellText'" }
$cell.Range.Text = "14+20=34"

$cell = $t.Cell(20, 5)
$cellText = $cell.Range.Text
if ($cellText.Substring(0, $cellText.Length - 2) -ne "39+1=40") { $mismatches++; Write-Output "MISMATCH at (20,5): expected '39+1=40' got '$cellText'" }
$cell.Range.Text = "84-55=29"

Write-Output ("Mismatches: " + $mismatches)
Write-Output "Done."